# The source workbook stores a tiny "is_active" flag table in Sheet1:
#   A1 = "is_active"   B1 = "false"
# This edit flips the flag to "True". The target value is the literal
# TEXT string "True" (not the Excel boolean TRUE), so we prefix the
# value with an apostrophe to force text entry and avoid Excel's
# automatic bool-literal coercion (same trick a human would use typing
# directly into the grid).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "'True"
